$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Splin")
$ws2 = $wb.Worksheets.Item("Iriao")

# ---------------------------------------------------------------------------
# Sheet "Splin": fix row 3, insert two new rows (4 and 5)
# ---------------------------------------------------------------------------

# Row 3: becomes another "Alexander Vasiliev / Russia" entry
$ws1.Range("B3").Value = 1
$ws1.Range("D3").Value = 7345934509
$ws1.Range("E3").Value = "Russia"
$ws1.Range("F3").Value = "Russian"
$ws1.Range("G3").Value = "Moscow"

# Row 4: new "Alexander Vasiliev / Russia" entry with a different birth date
$ws1.Range("A4").Value = "Alexander Vasiliev"
$ws1.Range("B4").Value = 1
$ws1.Range("C2").Copy($ws1.Range("C4"))
$ws1.Range("C4").Value = 21746
$ws1.Range("D4").Value = 7345934509
$ws1.Range("E4").Value = "Russia"
$ws1.Range("F4").Value = "Russian"
$ws1.Range("G4").Value = "Moscow"

# Row 5: the original Georgia/Tbilisi entry (shifted down from the old row 3)
$ws1.Range("A5").Value = "Alexander Vasiliev"
$ws1.Range("B5").Value = 1
$ws1.Range("C2").Copy($ws1.Range("C5"))
$ws1.Range("C5").Value = 25399
$ws1.Range("D5").Value = 7345934509
$ws1.Range("E5").Value = "Georgia"
$ws1.Range("F5").Value = "Georgian"
$ws1.Range("G5").Value = "Tbilisi"

# Extra (empty) column width tweak that shows up in the saved view state
$ws1.Columns.Item(8).ColumnWidth = 8

# ---------------------------------------------------------------------------
# Sheet "Iriao": append a near-duplicate row used to exercise the new
# uniqueness check
# ---------------------------------------------------------------------------

$ws2.Range("A4").Value = "Birdzina Muкia"
$ws2.Range("B4").Value = "авва"
$ws2.Range("C3").Copy($ws2.Range("C4"))
$ws2.Range("C4").Value = 29906
$ws2.Range("D4").Value = 6756453423
$ws2.Range("E4").Value = "Georgia"
$ws2.Range("F4").Value = "Georgian"
$ws2.Range("G4").Value = "Tbilisi"

# ---------------------------------------------------------------------------
# Selection / active-cell bookkeeping (matches the "error page" view change)
# ---------------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("F14").Select() | Out-Null

$ws2.Activate()
$ws2.Range("G8").Select() | Out-Null
